$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and volume-change (E) values from the latest GitHub Actions run.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.223.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2866"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07929"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6810"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.206.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.77%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007391"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.110.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.333"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.207"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.232"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.969"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09899"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.389"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.476"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.069"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04715"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7037"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01890"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.632"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.264"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4174"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "967.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.181"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.188"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05652"
$ws.Range("D51").Style = "Normal"
